{"js": "// 1. Name line: \"Name_____________________________________\" -> \"Name: Andrew Martin\"\nconst nameResults = context.document.body.search(\"Name_____________________________________\", { matchCase: true });\ncontext.load(nameResults);\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  const nameRange = nameResults.items[0];\n  nameRange.insertText(\"Name\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Re-search for the now-shortened \"Name\" run so we can append the rest right after it.\n  const nameResults2 = context.document.body.search(\"Name\", { matchCase: true });\n  context.load(nameResults2);\n  await context.sync();\n  nameResults2.items[0].insertText(\": Andrew Martin\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 2. Data table: fill in the empty \"L avg (m)\" column (index 3) with the average of\n//    the L1 (index 1) and L2 (index 2) columns for every data row.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const dataTable = tables.items[0];\n  dataTable.load(\"values,rowCount\");\n  await context.sync();\n\n  const values = dataTable.values;\n  for (let row = 1; row < values.length; row++) {\n    const l1 = parseFloat(values[row][1]);\n    const l2 = parseFloat(values[row][2]);\n    const existingAvg = (values[row][3] || \"\").trim();\n    if (!isNaN(l1) && !isNaN(l2) && existingAvg === \"\") {\n      let avg = (l1 + l2) / 2;\n      // Trim floating point noise (e.g. 0.28 + 0.29 -> 0.285 exactly, not 0.28500000000000003)\n      avg = Math.round(avg * 1e6) / 1e6;\n      const avgText = String(avg);\n\n      const cell = dataTable.getCell(row, 3);\n      const paragraph = cell.body.paragraphs.getFirst();\n      const range = paragraph.getRange(\"Whole\");\n      range.insertText(avgText, Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Name line: \"Name_____________________________________\" -> \"Name: Andrew Martin\"\n$find = $d.Content.Find\n$find.Text = \"Name_____________________________________\"\n$found = $find.Execute()\nif ($found) {\n    $rng = $find.Parent\n    $rng.Text = \"Name\"\n    $rng.Collapse(0)\n    $rng.InsertAfter(\": Andrew Martin\")\n}\n\n# 2. Data table: fill in the empty \"L avg (m)\" column (column 4) with the average of\n#    the L1 (column 2) and L2 (column 3) columns for every data row.\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n\nfor ($i = 2; $i -le $rowCount; $i++) {\n    $l1Cell = $tbl.Cell($i, 2)\n    $l2Cell = $tbl.Cell($i, 3)\n    $avgCell = $tbl.Cell($i, 4)\n\n    $l1Text = $l1Cell.Range.Text.Trim([char]7, [char]13, [char]10)\n    $l2Text = $l2Cell.Range.Text.Trim([char]7, [char]13, [char]10)\n    $avgText = $avgCell.Range.Text.Trim([char]7, [char]13, [char]10)\n\n    if ($avgText -eq \"\" -and $l1Text -ne \"\" -and $l2Text -ne \"\") {\n        $l1 = [double]$l1Text\n        $l2 = [double]$l2Text\n        $avg = [Math]::Round((($l1 + $l2) / 2), 6)\n\n        $avgPara = $avgCell.Range.Paragraphs.Item(1)\n        $avgRng = $avgPara.Range\n        $avgRng.Text = $avg.ToString()\n    }\n}\n"}
